$d = $word.ActiveDocument

# Locate the paragraph that holds "URL:" followed by the hyperlink
# (it is the paragraph whose text starts with "URL:").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("URL:")) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the URL paragraph"
}

# Replace the whole paragraph (including its end-of-paragraph mark) with
# two paragraphs: the first one keeps "URL: <hyperlink>" (now carrying
# es-ES language tags + hyperlink-style paragraph mark formatting), and
# the second one is a brand new paragraph with a hyperlink-styled
# "CambiarURL" run; the _GoBack bookmark now wraps that new run instead.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:rPr>' + `
      '<w:rStyle w:val="Hipervnculo"/>' + `
      '<w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/>' + `
      '<w:color w:val="1155CC"/>' + `
      '<w:sz w:val="19"/>' + `
      '<w:szCs w:val="19"/>' + `
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
      '<w:lang w:val="es-ES"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:b/>' + `
      '<w:sz w:val="28"/>' + `
      '<w:szCs w:val="28"/>' + `
      '<w:lang w:val="es-ES"/>' + `
    '</w:rPr>' + `
    '<w:t>URL:</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:lang w:val="es-ES"/>' + `
    '</w:rPr>' + `
    '<w:t xml:space="preserve"> </w:t>' + `
  '</w:r>' + `
  '<w:hyperlink r:id="rId5" w:tgtFrame="_blank" w:history="1">' + `
    '<w:r>' + `
      '<w:rPr>' + `
        '<w:rStyle w:val="Hipervnculo"/>' + `
        '<w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/>' + `
        '<w:color w:val="1155CC"/>' + `
        '<w:sz w:val="19"/>' + `
        '<w:szCs w:val="19"/>' + `
        '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
        '<w:lang w:val="es-ES"/>' + `
      '</w:rPr>' + `
      '<w:t>http://www.hlp.somee.com/Inicio.aspx</w:t>' + `
    '</w:r>' + `
  '</w:hyperlink>' + `
'</w:p>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:rPr>' + `
      '<w:lang w:val="es-ES"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rStyle w:val="Hipervnculo"/>' + `
      '<w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/>' + `
      '<w:color w:val="1155CC"/>' + `
      '<w:sz w:val="19"/>' + `
      '<w:szCs w:val="19"/>' + `
      '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' + `
      '<w:lang w:val="es-ES"/>' + `
    '</w:rPr>' + `
    '<w:t>CambiarURL</w:t>' + `
  '</w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
'</w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)

# InsertXML does not preserve <w:rStyle> references (the rest of the
# direct character formatting - fonts/color/size/shading - survives the
# round trip fine), so re-apply the "Hipervnculo" character style to the
# two runs that need it via a scoped Find/Replace (this is the one COM
# path that reliably emits <w:rStyle w:val="Hipervnculo"/>).
function Set-HyperlinkCharStyle($searchRange, $text) {
    $searchRange.Find.ClearFormatting()
    $searchRange.Find.Replacement.ClearFormatting()
    $searchRange.Find.Replacement.Style = "Hipervnculo"
    $searchRange.Find.Execute($text, $false, $false, $false, $false, $false, `
                               $true, 1, $false, $text, 2)
}

Set-HyperlinkCharStyle $d.Content "http://www.hlp.somee.com/Inicio.aspx"
Set-HyperlinkCharStyle $d.Content "CambiarURL"
